$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update end time for the row 7 entry (Flutter.dev introduction day) from 15 to 16,
# which cascades through the shared "time worked" formula and the summary totals.
$ws.Range("C7").Value = 16

$wb.Application.Calculate()
